$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(28, 8).Value = 5833.615
$ws.Cells.Item(28, 9).Value = 1373.8
$ws.Cells.Item(28, 11).Value = 1373.8
$ws.Cells.Item(28, 13).Value = -888.8
$ws.Cells.Item(33, 8).Value = 5555778.5
$ws.Cells.Item(33, 9).Value = 7143098
$ws.Cells.Item(33, 10).Value = 159.25
$ws.Cells.Item(33, 11).Value = 7143098
$ws.Cells.Item(33, 12).Value = 159.25
$ws.Cells.Item(33, 13).Value = -7142869
$ws.Cells.Item(33, 14).Value = -617.25
$ws.Cells.Item(62, 8).Value = 3014.4043
$ws.Cells.Item(62, 9).Value = 2728.756
$ws.Cells.Item(62, 10).Value = 4966.3335
$ws.Cells.Item(62, 11).Value = 2728.756
$ws.Cells.Item(62, 12).Value = 4966.3335
$ws.Cells.Item(62, 13).Value = -2104.756
$ws.Cells.Item(62, 14).Value = -6214.3335
$ws.Cells.Item(65, 8).Value = 3014.4043
$ws.Cells.Item(65, 9).Value = 2728.756
$ws.Cells.Item(65, 10).Value = 4966.3335
$ws.Cells.Item(65, 11).Value = 13643.78
$ws.Cells.Item(65, 12).Value = 24831.6675
$ws.Cells.Item(65, 13).Value = -10523.78
$ws.Cells.Item(65, 14).Value = -31071.6675
$ws.Cells.Item(70, 8).Value = 4531.4165
$ws.Cells.Item(70, 10).Value = 4630.778
$ws.Cells.Item(70, 12).Value = 13892.334
$ws.Cells.Item(70, 14).Value = -14432.334
$ws.Cells.Item(73, 8).Value = 4531.4165
$ws.Cells.Item(73, 10).Value = 4630.778
$ws.Cells.Item(73, 12).Value = 13892.334
$ws.Cells.Item(73, 14).Value = -15764.334
$ws.Cells.Item(76, 8).Value = 5125
$ws.Cells.Item(76, 9).Value = 4666.6665
$ws.Cells.Item(76, 10).Value = 6500
$ws.Cells.Item(76, 11).Value = 4666.6665
$ws.Cells.Item(76, 12).Value = 6500
$ws.Cells.Item(76, 13).Value = -4351.6665
$ws.Cells.Item(76, 14).Value = -7130
$ws.Cells.Item(79, 8).Value = 5125
$ws.Cells.Item(79, 9).Value = 4666.6665
$ws.Cells.Item(79, 10).Value = 6500
$ws.Cells.Item(79, 11).Value = 4666.6665
$ws.Cells.Item(79, 12).Value = 6500
$ws.Cells.Item(79, 13).Value = -3574.6665
$ws.Cells.Item(79, 14).Value = -8684
$ws.Cells.Item(80, 8).Value = 1490.4166
$ws.Cells.Item(80, 9).Value = 2066.6667
$ws.Cells.Item(80, 10).Value = 1298.3334
$ws.Cells.Item(80, 11).Value = 6200.000100000001
$ws.Cells.Item(80, 12).Value = 3895.0002
$ws.Cells.Item(80, 13).Value = -5202.000100000001
$ws.Cells.Item(80, 14).Value = -5891.0002
$ws.Cells.Item(83, 8).Value = 1490.4166
$ws.Cells.Item(83, 9).Value = 2066.6667
$ws.Cells.Item(83, 10).Value = 1298.3334
$ws.Cells.Item(83, 11).Value = 18600.0003
$ws.Cells.Item(83, 12).Value = 11685.0006
$ws.Cells.Item(83, 13).Value = -13608.0003
$ws.Cells.Item(83, 14).Value = -21669.0006
$ws.Cells.Item(86, 8).Value = 2549.6667
$ws.Cells.Item(86, 9).Value = 2699.5
$ws.Cells.Item(86, 11).Value = 2699.5
$ws.Cells.Item(86, 13).Value = -1576.5
$ws.Cells.Item(89, 8).Value = 2549.6667
$ws.Cells.Item(89, 9).Value = 2699.5
$ws.Cells.Item(89, 11).Value = 13497.5
$ws.Cells.Item(89, 13).Value = -7881.5
$ws.Cells.Item(107, 8).Value = 16131058
$ws.Cells.Item(107, 9).Value = 20834034
$ws.Cells.Item(107, 10).Value = 6570.5713
$ws.Cells.Item(107, 11).Value = 20834034
$ws.Cells.Item(107, 12).Value = 6570.5713
$ws.Cells.Item(107, 13).Value = -20832114
$ws.Cells.Item(107, 14).Value = -10410.5713
$ws.Cells.Item(132, 8).Value = 2389.3076
$ws.Cells.Item(132, 9).Value = 1978.125
$ws.Cells.Item(132, 11).Value = 5934.375
$ws.Cells.Item(132, 13).Value = -3404.375
$ws.Cells.Item(134, 8).Value = 103593.336
$ws.Cells.Item(134, 9).Value = 90000
$ws.Cells.Item(134, 10).Value = 110390
$ws.Cells.Item(134, 11).Value = 90000
$ws.Cells.Item(134, 12).Value = 110390
$ws.Cells.Item(134, 13).Value = -84930
$ws.Cells.Item(134, 14).Value = -120530
$ws.Cells.Item(137, 8).Value = 7583888.5
$ws.Cells.Item(137, 9).Value = 16670594
$ws.Cells.Item(137, 10).Value = 11634.223
$ws.Cells.Item(137, 11).Value = 50011782
$ws.Cells.Item(137, 12).Value = 34902.669
$ws.Cells.Item(137, 13).Value = -50009232
$ws.Cells.Item(137, 14).Value = -40002.669
$ws.Cells.Item(138, 8).Value = 4248.028
$ws.Cells.Item(138, 9).Value = 2413.4707
$ws.Cells.Item(138, 10).Value = 5889.4736
$ws.Cells.Item(138, 11).Value = 7240.4121
$ws.Cells.Item(138, 12).Value = 17668.4208
$ws.Cells.Item(138, 13).Value = -2100.4121
$ws.Cells.Item(138, 14).Value = -27948.4208
$ws.Cells.Item(141, 8).Value = 832.8
$ws.Cells.Item(141, 9).Value = 832.8
$ws.Cells.Item(141, 11).Value = 2498.4
$ws.Cells.Item(141, 13).Value = 2681.6
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(2, 8).Value = 2675244.5
$ws.Cells.Item(2, 9).Value = 2842244.5
$ws.Cells.Item(2, 11).Value = 2842244.5
$ws.Cells.Item(2, 13).Value = -2842131.5
$ws.Cells.Item(32, 8).Value = 12582.292
$ws.Cells.Item(32, 9).Value = 9232.471
$ws.Cells.Item(32, 11).Value = 9232.471
$ws.Cells.Item(32, 13).Value = -8945.471
$ws.Cells.Item(34, 8).Value = 181397.8
$ws.Cells.Item(34, 9).Value = 101747.25
$ws.Cells.Item(34, 11).Value = 101747.25
$ws.Cells.Item(34, 13).Value = -101476.25
$ws.Cells.Item(38, 8).Value = 0
$ws.Cells.Item(38, 9).Value = 0
$ws.Cells.Item(38, 11).Value = 0
$ws.Cells.Item(38, 13).Value = ""
$ws.Cells.Item(39, 8).Value = 5906.2
$ws.Cells.Item(39, 9).Value = 2882.75
$ws.Cells.Item(39, 10).Value = 18000
$ws.Cells.Item(39, 11).Value = 2882.75
$ws.Cells.Item(39, 12).Value = 18000
$ws.Cells.Item(39, 13).Value = -2362.75
$ws.Cells.Item(39, 14).Value = -19040
$ws.Cells.Item(40, 8).Value = 20031
$ws.Cells.Item(40, 10).Value = 20031
$ws.Cells.Item(40, 12).Value = 20031
$ws.Cells.Item(40, 14).Value = -20383
$ws.Cells.Item(42, 8).Value = 19030.5
$ws.Cells.Item(42, 9).Value = 19030
$ws.Cells.Item(42, 11).Value = 19030
$ws.Cells.Item(42, 13).Value = -18544
$ws.Cells.Item(45, 8).Value = 1322.2
$ws.Cells.Item(45, 9).Value = 1215.25
$ws.Cells.Item(45, 11).Value = 1215.25
$ws.Cells.Item(45, 13).Value = -838.25
$ws.Cells.Item(61, 8).Value = 5289.875
$ws.Cells.Item(61, 9).Value = 3402.7144
$ws.Cells.Item(61, 11).Value = 3402.7144
$ws.Cells.Item(61, 13).Value = -3190.7144
$ws.Cells.Item(74, 8).Value = 5913
$ws.Cells.Item(74, 9).Value = 2310.75
$ws.Cells.Item(74, 11).Value = 2310.75
$ws.Cells.Item(74, 13).Value = -1436.75
$ws.Cells.Item(77, 8).Value = 5913
$ws.Cells.Item(77, 9).Value = 2310.75
$ws.Cells.Item(77, 11).Value = 11553.75
$ws.Cells.Item(77, 13).Value = -7185.75
$ws.Cells.Item(116, 8).Value = 2675244.5
$ws.Cells.Item(116, 9).Value = 2842244.5
$ws.Cells.Item(116, 11).Value = 2842244.5
$ws.Cells.Item(116, 13).Value = -2839950.5
$ws.Cells.Item(122, 8).Value = 3055.5144
$ws.Cells.Item(122, 9).Value = 3013.9062
$ws.Cells.Item(122, 10).Value = 3499.3333
$ws.Cells.Item(122, 11).Value = 9041.7186
$ws.Cells.Item(122, 12).Value = 10497.9999
$ws.Cells.Item(122, 13).Value = -6591.7186
$ws.Cells.Item(122, 14).Value = -15397.9999
$ws.Cells.Item(132, 8).Value = 5046.5776
$ws.Cells.Item(132, 9).Value = 2366.6365
$ws.Cells.Item(132, 10).Value = 12416.417
$ws.Cells.Item(132, 11).Value = 7099.9095
$ws.Cells.Item(132, 12).Value = 37249.251
$ws.Cells.Item(132, 13).Value = -4569.9095
$ws.Cells.Item(132, 14).Value = -42309.251
$ws.Cells.Item(136, 8).Value = 5289.875
$ws.Cells.Item(136, 9).Value = 3402.7144
$ws.Cells.Item(136, 11).Value = 10208.1432
$ws.Cells.Item(136, 13).Value = -7658.143199999999
$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(3, 8).Value = 2675244.5
$ws.Cells.Item(3, 9).Value = 2842244.5
$ws.Cells.Item(3, 11).Value = 2842244.5
$ws.Cells.Item(3, 13).Value = -2842130.5
$ws.Cells.Item(20, 8).Value = 3219.853
$ws.Cells.Item(20, 9).Value = 3198.2778
$ws.Cells.Item(20, 10).Value = 3244.125
$ws.Cells.Item(20, 11).Value = 3198.2778
$ws.Cells.Item(20, 12).Value = 3244.125
$ws.Cells.Item(20, 13).Value = -2951.2778
$ws.Cells.Item(20, 14).Value = -3738.125
$ws.Cells.Item(24, 8).Value = 800
$ws.Cells.Item(24, 9).Value = 800
$ws.Cells.Item(24, 10).Value = 0
$ws.Cells.Item(24, 11).Value = 800
$ws.Cells.Item(24, 12).Value = 0
$ws.Cells.Item(24, 13).Value = -565
$ws.Cells.Item(24, 14).Value = ""
$ws.Cells.Item(25, 8).Value = 6630.231
$ws.Cells.Item(25, 9).Value = 399.25
$ws.Cells.Item(25, 11).Value = 399.25
$ws.Cells.Item(25, 13).Value = -164.25
$ws.Cells.Item(40, 8).Value = 20440
$ws.Cells.Item(40, 10).Value = 20440
$ws.Cells.Item(40, 12).Value = 20440
$ws.Cells.Item(40, 14).Value = -20970
$ws.Cells.Item(80, 8).Value = 50533.3
$ws.Cells.Item(80, 10).Value = 83738.086
$ws.Cells.Item(80, 12).Value = 83738.086
$ws.Cells.Item(80, 14).Value = -85734.086
$ws.Cells.Item(83, 8).Value = 50533.3
$ws.Cells.Item(83, 10).Value = 83738.086
$ws.Cells.Item(83, 12).Value = 418690.43
$ws.Cells.Item(83, 14).Value = -428674.43
$ws.Cells.Item(86, 8).Value = 2003.3182
$ws.Cells.Item(86, 9).Value = 1556.4736
$ws.Cells.Item(86, 11).Value = 1556.4736
$ws.Cells.Item(86, 13).Value = -433.4736
$ws.Cells.Item(89, 8).Value = 2003.3182
$ws.Cells.Item(89, 9).Value = 1556.4736
$ws.Cells.Item(89, 11).Value = 7782.368
$ws.Cells.Item(89, 13).Value = -2166.368
$ws.Cells.Item(94, 8).Value = 602.40625
$ws.Cells.Item(94, 9).Value = 471.6207
$ws.Cells.Item(94, 10).Value = 1866.6666
$ws.Cells.Item(94, 11).Value = 471.6207
$ws.Cells.Item(94, 12).Value = 1866.6666
$ws.Cells.Item(94, 13).Value = -20.6207
$ws.Cells.Item(94, 14).Value = -2768.6666
$ws.Cells.Item(105, 8).Value = 83335570
$ws.Cells.Item(105, 9).Value = 100002020
$ws.Cells.Item(105, 10).Value = 3305.5
$ws.Cells.Item(105, 11).Value = 100002020
$ws.Cells.Item(105, 12).Value = 3305.5
$ws.Cells.Item(105, 13).Value = -100000273
$ws.Cells.Item(105, 14).Value = -6799.5
$ws.Cells.Item(107, 8).Value = 1710.4783
$ws.Cells.Item(107, 9).Value = 1657.421
$ws.Cells.Item(107, 10).Value = 1962.5
$ws.Cells.Item(107, 11).Value = 1657.421
$ws.Cells.Item(107, 12).Value = 1962.5
$ws.Cells.Item(107, 13).Value = 262.579
$ws.Cells.Item(107, 14).Value = -5802.5
$ws.Cells.Item(132, 8).Value = 0
$ws.Cells.Item(132, 10).Value = 0
$ws.Cells.Item(132, 12).Value = 0
$ws.Cells.Item(132, 14).Value = ""
$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(92, 8).Value = 31000
$ws.Cells.Item(92, 10).Value = 31000
$ws.Cells.Item(92, 12).Value = 31000
$ws.Cells.Item(92, 14).Value = -35992
$ws.Cells.Item(99, 8).Value = 1554.6
$ws.Cells.Item(99, 9).Value = 1658
$ws.Cells.Item(99, 10).Value = 1399.5
$ws.Cells.Item(99, 11).Value = 1658
$ws.Cells.Item(99, 12).Value = 1399.5
$ws.Cells.Item(99, 13).Value = -160
$ws.Cells.Item(99, 14).Value = -4395.5
$ws.Cells.Item(126, 8).Value = 1554.6
$ws.Cells.Item(126, 9).Value = 1658
$ws.Cells.Item(126, 10).Value = 1399.5
$ws.Cells.Item(126, 11).Value = 4974
$ws.Cells.Item(126, 12).Value = 4198.5
$ws.Cells.Item(126, 13).Value = -2504
$ws.Cells.Item(126, 14).Value = -9138.5
$ws = $wb.Worksheets.Item(5)
$ws.Cells.Item(2, 8).Value = 2432.3667
$ws.Cells.Item(2, 9).Value = 77.73333
$ws.Cells.Item(2, 10).Value = 4787
$ws.Cells.Item(2, 11).Value = 466.39998
$ws.Cells.Item(2, 12).Value = 28722
$ws.Cells.Item(2, 13).Value = -353.39998
$ws.Cells.Item(2, 14).Value = -28948
$ws.Cells.Item(5, 8).Value = 336.80768
$ws.Cells.Item(5, 10).Value = 399.35294
$ws.Cells.Item(5, 12).Value = 1198.05882
$ws.Cells.Item(5, 14).Value = -1422.05882
$ws.Cells.Item(13, 8).Value = 368.85715
$ws.Cells.Item(13, 9).Value = 277.33334
$ws.Cells.Item(13, 10).Value = 437.5
$ws.Cells.Item(13, 11).Value = 832.0000200000001
$ws.Cells.Item(13, 12).Value = 1312.5
$ws.Cells.Item(13, 13).Value = -664.0000200000001
$ws.Cells.Item(13, 14).Value = -1648.5
$ws.Cells.Item(38, 8).Value = 36.57895
$ws.Cells.Item(38, 9).Value = 10.714286
$ws.Cells.Item(38, 10).Value = 51.666668
$ws.Cells.Item(38, 11).Value = 32.142858
$ws.Cells.Item(38, 12).Value = 155.000004
$ws.Cells.Item(38, 13).Value = 314.857142
$ws.Cells.Item(38, 14).Value = -849.000004
$ws.Cells.Item(56, 8).Value = 6555.25
$ws.Cells.Item(56, 9).Value = 6555.25
$ws.Cells.Item(56, 11).Value = 6555.25
$ws.Cells.Item(56, 13).Value = -6025.25
$ws.Cells.Item(100, 8).Value = 988.46155
$ws.Cells.Item(100, 10).Value = 988.46155
$ws.Cells.Item(100, 12).Value = 2965.38465
$ws.Cells.Item(100, 14).Value = -4587.38465
$ws.Cells.Item(108, 8).Value = 3170
$ws.Cells.Item(108, 9).Value = 3170
$ws.Cells.Item(108, 11).Value = 9510
$ws.Cells.Item(108, 13).Value = -6630
$ws.Cells.Item(113, 8).Value = 3023.8
$ws.Cells.Item(113, 10).Value = 6585
$ws.Cells.Item(113, 12).Value = 19755
$ws.Cells.Item(113, 14).Value = -24095
$ws.Cells.Item(131, 8).Value = 4608.3184
$ws.Cells.Item(131, 9).Value = 932.6667
$ws.Cells.Item(131, 10).Value = 5188.684
$ws.Cells.Item(131, 11).Value = 2798.0001
$ws.Cells.Item(131, 12).Value = 15566.052
$ws.Cells.Item(131, 13).Value = 2241.9999
$ws.Cells.Item(131, 14).Value = -25646.052
$ws.Cells.Item(132, 8).Value = 1692.6
$ws.Cells.Item(132, 9).Value = 1164.5
$ws.Cells.Item(132, 10).Value = 2044.6666
$ws.Cells.Item(132, 11).Value = 10480.5
$ws.Cells.Item(132, 12).Value = 18401.9994
$ws.Cells.Item(132, 13).Value = -7950.5
$ws.Cells.Item(132, 14).Value = -23461.9994
$ws.Cells.Item(135, 8).Value = 336.80768
$ws.Cells.Item(135, 10).Value = 399.35294
$ws.Cells.Item(135, 12).Value = 3594.17646
$ws.Cells.Item(135, 14).Value = -8664.176459999999
$ws = $wb.Worksheets.Item(6)
$ws.Cells.Item(23, 8).Value = 849.75
$ws.Cells.Item(23, 10).Value = 799
$ws.Cells.Item(23, 12).Value = 799
$ws.Cells.Item(23, 14).Value = -1245
$ws.Cells.Item(80, 8).Value = 2586.2222
$ws.Cells.Item(80, 9).Value = 2355.1428
$ws.Cells.Item(80, 10).Value = 3395
$ws.Cells.Item(80, 11).Value = 2355.1428
$ws.Cells.Item(80, 12).Value = 3395
$ws.Cells.Item(80, 13).Value = -1357.1428
$ws.Cells.Item(80, 14).Value = -5391
$ws.Cells.Item(83, 8).Value = 2586.2222
$ws.Cells.Item(83, 9).Value = 2355.1428
$ws.Cells.Item(83, 10).Value = 3395
$ws.Cells.Item(83, 11).Value = 11775.714
$ws.Cells.Item(83, 12).Value = 16975
$ws.Cells.Item(83, 13).Value = -6783.714
$ws.Cells.Item(83, 14).Value = -26959
$ws.Cells.Item(102, 8).Value = 1773.5
$ws.Cells.Item(102, 9).Value = 1900.6666
$ws.Cells.Item(102, 11).Value = 1900.6666
$ws.Cells.Item(102, 13).Value = -278.6666
$ws.Cells.Item(113, 8).Value = 92698.60000000001
$ws.Cells.Item(113, 9).Value = 6750
$ws.Cells.Item(113, 10).Value = 149997.67
$ws.Cells.Item(113, 11).Value = 6750
$ws.Cells.Item(113, 12).Value = 149997.67
$ws.Cells.Item(113, 13).Value = -4580
$ws.Cells.Item(113, 14).Value = -154337.67
$ws.Cells.Item(122, 8).Value = 5491.381
$ws.Cells.Item(122, 9).Value = 5858.8945
$ws.Cells.Item(122, 10).Value = 2000
$ws.Cells.Item(122, 11).Value = 17576.6835
$ws.Cells.Item(122, 12).Value = 6000
$ws.Cells.Item(122, 13).Value = -15126.6835
$ws.Cells.Item(122, 14).Value = -10900
$ws.Cells.Item(126, 8).Value = 3200.3
$ws.Cells.Item(126, 9).Value = 1913
$ws.Cells.Item(126, 11).Value = 5739
$ws.Cells.Item(126, 13).Value = -3269
$ws.Cells.Item(132, 8).Value = 7051.1904
$ws.Cells.Item(132, 9).Value = 5304.645
$ws.Cells.Item(132, 11).Value = 15913.935
$ws.Cells.Item(132, 13).Value = -13383.935
$ws.Cells.Item(140, 8).Value = 77854.25
$ws.Cells.Item(140, 10).Value = 94999.5
$ws.Cells.Item(140, 12).Value = 94999.5
$ws.Cells.Item(140, 14).Value = -105359.5
$ws = $wb.Worksheets.Item(7)
$ws.Cells.Item(22, 8).Value = 3079.0715
$ws.Cells.Item(22, 10).Value = 4221.0527
$ws.Cells.Item(22, 12).Value = 4221.0527
$ws.Cells.Item(22, 14).Value = -4811.0527
$ws.Cells.Item(27, 8).Value = 3079.0715
$ws.Cells.Item(27, 10).Value = 4221.0527
$ws.Cells.Item(27, 12).Value = 4221.0527
$ws.Cells.Item(27, 14).Value = -4435.0527
$ws.Cells.Item(33, 8).Value = 53000
$ws.Cells.Item(33, 10).Value = 53000
$ws.Cells.Item(33, 12).Value = 53000
$ws.Cells.Item(33, 14).Value = -53580
$ws.Cells.Item(40, 8).Value = 5826.4736
$ws.Cells.Item(40, 9).Value = 5483.5
$ws.Cells.Item(40, 11).Value = 5483.5
$ws.Cells.Item(40, 13).Value = -5347.5
$ws.Cells.Item(68, 8).Value = 2580
$ws.Cells.Item(68, 9).Value = 2596
$ws.Cells.Item(68, 10).Value = 2500
$ws.Cells.Item(68, 11).Value = 2596
$ws.Cells.Item(68, 12).Value = 2500
$ws.Cells.Item(68, 13).Value = -1847
$ws.Cells.Item(68, 14).Value = -3998
$ws.Cells.Item(71, 8).Value = 2580
$ws.Cells.Item(71, 9).Value = 2596
$ws.Cells.Item(71, 10).Value = 2500
$ws.Cells.Item(71, 11).Value = 12980
$ws.Cells.Item(71, 12).Value = 12500
$ws.Cells.Item(71, 13).Value = -9236
$ws.Cells.Item(71, 14).Value = -19988
$ws.Cells.Item(82, 8).Value = 1374
$ws.Cells.Item(82, 9).Value = 1011
$ws.Cells.Item(82, 10).Value = 2100
$ws.Cells.Item(82, 11).Value = 1011
$ws.Cells.Item(82, 12).Value = 2100
$ws.Cells.Item(82, 13).Value = -650
$ws.Cells.Item(82, 14).Value = -2822
$ws.Cells.Item(85, 8).Value = 1374
$ws.Cells.Item(85, 9).Value = 1011
$ws.Cells.Item(85, 10).Value = 2100
$ws.Cells.Item(85, 11).Value = 1011
$ws.Cells.Item(85, 12).Value = 2100
$ws.Cells.Item(85, 13).Value = 237
$ws.Cells.Item(85, 14).Value = -4596
$ws.Cells.Item(100, 8).Value = 7816286
$ws.Cells.Item(100, 9).Value = 13892031
$ws.Cells.Item(100, 11).Value = 13892031
$ws.Cells.Item(100, 13).Value = -13891490
$ws.Cells.Item(115, 8).Value = 62500
$ws.Cells.Item(115, 10).Value = 62500
$ws.Cells.Item(115, 12).Value = 62500
$ws.Cells.Item(115, 14).Value = -64850
$ws.Cells.Item(118, 8).Value = 48000
$ws.Cells.Item(118, 10).Value = 48000
$ws.Cells.Item(118, 12).Value = 48000
$ws.Cells.Item(118, 14).Value = -51314
$ws.Cells.Item(122, 8).Value = 6334.3335
$ws.Cells.Item(122, 9).Value = 8002
$ws.Cells.Item(122, 10).Value = 2999
$ws.Cells.Item(122, 11).Value = 24006
$ws.Cells.Item(122, 12).Value = 8997
$ws.Cells.Item(122, 13).Value = -21556
$ws.Cells.Item(122, 14).Value = -13897
$ws.Cells.Item(132, 8).Value = 5766.6113
$ws.Cells.Item(132, 9).Value = 4475
$ws.Cells.Item(132, 11).Value = 13425
$ws.Cells.Item(132, 13).Value = -10895
$ws.Cells.Item(136, 8).Value = 9208.333000000001
$ws.Cells.Item(136, 9).Value = 1500
$ws.Cells.Item(136, 10).Value = 9909.091
$ws.Cells.Item(136, 11).Value = 4500
$ws.Cells.Item(136, 12).Value = 29727.273
$ws.Cells.Item(136, 13).Value = -1950
$ws.Cells.Item(136, 14).Value = -34827.273
$ws = $wb.Worksheets.Item(8)
$ws.Cells.Item(15, 8).Value = 0
$ws.Cells.Item(15, 9).Value = 0
$ws.Cells.Item(15, 11).Value = 0
$ws.Cells.Item(15, 13).Value = ""
$ws.Cells.Item(43, 8).Value = 67577.8
$ws.Cells.Item(43, 9).Value = 48945
$ws.Cells.Item(43, 10).Value = 79999.664
$ws.Cells.Item(43, 11).Value = 48945
$ws.Cells.Item(43, 12).Value = 79999.664
$ws.Cells.Item(43, 13).Value = -48796
$ws.Cells.Item(43, 14).Value = -80297.664
$ws.Cells.Item(103, 8).Value = 53766.855
$ws.Cells.Item(103, 10).Value = 53766.855
$ws.Cells.Item(103, 12).Value = 53766.855
$ws.Cells.Item(103, 14).Value = -56110.855
$ws.Cells.Item(107, 8).Value = 2115.718
$ws.Cells.Item(107, 9).Value = 2160.2424
$ws.Cells.Item(107, 10).Value = 1870.8334
$ws.Cells.Item(107, 11).Value = 6480.7272
$ws.Cells.Item(107, 12).Value = 5612.5002
$ws.Cells.Item(107, 13).Value = -4560.7272
$ws.Cells.Item(107, 14).Value = -9452.5002
$ws.Cells.Item(122, 8).Value = 3442.4546
$ws.Cells.Item(122, 9).Value = 3116.8235
$ws.Cells.Item(122, 10).Value = 4549.6
$ws.Cells.Item(122, 11).Value = 9350.470499999999
$ws.Cells.Item(122, 12).Value = 13648.8
$ws.Cells.Item(122, 13).Value = -6900.470499999999
$ws.Cells.Item(122, 14).Value = -18548.8
$ws.Cells.Item(126, 8).Value = 3794.0435
$ws.Cells.Item(126, 9).Value = 4083.9524
$ws.Cells.Item(126, 10).Value = 750
$ws.Cells.Item(126, 11).Value = 12251.8572
$ws.Cells.Item(126, 12).Value = 2250
$ws.Cells.Item(126, 13).Value = -9781.8572
$ws.Cells.Item(126, 14).Value = -7190
$ws.Cells.Item(132, 8).Value = 5018.648
$ws.Cells.Item(132, 9).Value = 3113.9688
$ws.Cells.Item(132, 10).Value = 7789.091
$ws.Cells.Item(132, 11).Value = 9341.9064
$ws.Cells.Item(132, 12).Value = 23367.273
$ws.Cells.Item(132, 13).Value = -6811.9064
$ws.Cells.Item(132, 14).Value = -28427.273
$ws.Cells.Item(141, 8).Value = 120000
$ws.Cells.Item(141, 10).Value = 120000
$ws.Cells.Item(141, 12).Value = 120000
$ws.Cells.Item(141, 14).Value = -130360
